$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -21.988
    14 = -21.77059999999999
    16 = -21.63439999999999
    21 = -20.02659999999998
    23 = -20.55909999999998
    25 = -21.68069999999998
    26 = -21.04699999999997
    29 = -20.74099999999997
    40 = -19.73239999999999
    53 = -21.821
    57 = -22.20860000000001
    59 = -22.26549999999999
    65 = -21.80649999999998
    69 = -21.6255
    79 = -20.3359
    83 = -21.8233
    91 = -21.46840000000002
    93 = -20.71299999999998
    100 = -21.7327
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
